# Weekly update: prepend a new price record for "Bruselas (repollito)" at
# Vega Modelo de Temuco, pushing the existing rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 32 (row 1 is the header, data starts at row 2),
# shifting the existing rows 32:41 down to 33:42.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row with this week's data point.
$ws.Cells.Item(32, 1).Value = 10
$ws.Cells.Item(32, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(32, 3).Value = "La Araucanía"
$ws.Cells.Item(32, 4).Value = 44438
$ws.Cells.Item(32, 5).Value = 9
$ws.Cells.Item(32, 6).Value = 100112035
$ws.Cells.Item(32, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 50
$ws.Cells.Item(32, 11).Value = 25000
$ws.Cells.Item(32, 12).Value = 25000
$ws.Cells.Item(32, 13).Value = 25000
$ws.Cells.Item(32, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(32, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(32, 16).Value = 2500
$ws.Cells.Item(32, 17).Value = 10
$ws.Cells.Item(32, 18).Value = "Hortaliza"
